$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) row values for rows 2, 3, 4 across the
# columns that participate in the cyclic rotation described by the diff:
# new row2 = old row3, new row3 = old row4, new row4 = old row2.
$cols = @("A","B","D","E","F","G","H","Q","R")

$row2 = @{}
$row3 = @{}
$row4 = @{}

foreach ($col in $cols) {
    $row2[$col] = $ws.Range($col + "2").Value2
    $row3[$col] = $ws.Range($col + "3").Value2
    $row4[$col] = $ws.Range($col + "4").Value2
}

foreach ($col in $cols) {
    $ws.Range($col + "2").Value2 = $row3[$col]
    $ws.Range($col + "3").Value2 = $row4[$col]
    $ws.Range($col + "4").Value2 = $row2[$col]
}
